# Regenerate merged AHB files
# - Rename the "_old"/"_new" header-suffix columns to "_FV2404"/"_FV2410"
# - Turn the data range into an Excel Table ("Table1")
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) -----------------------------------------
# Columns A:J -> "..._old" becomes "..._FV2404"
# Column  K   -> "diff" stays as-is
# Columns L:U -> "..._new" becomes "..._FV2410"
$oldHeaders = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)
$newHeaders = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# --- 2. Convert the used range into a Table --------------------------------
$tableRange = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
